# The only content-level change in the target diff is the text of the
# (sole) shared string used by B2:B16: "dummy１" -> "dummy２".
# (The other hunks in the diff -- absPath, revisionPtr/documentId,
# workbookView window geometry, and the sheetView <selection> element --
# are Excel-session/environment metadata that isn't driven by the content
# object model, so we leave those alone and just rewrite the cell text.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("B2:B16")
$rng.Value = "dummy２"
